$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = 5
$ws.Range("D1").Value = 100
$ws.Range("C2").Value = 10
$ws.Range("C3").Value = 30
